$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot (GitHub Actions cron update).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Price cells that look numeric get forced to Text format first so Excel
# doesn't silently convert strings like "6.24" into floating point numbers
# (values that already contain two dots, e.g. "67.280.04", can never be
# parsed as a number so they are left alone).

$ws.Cells.Item(2, 4).Value = "67.280.04"
$ws.Cells.Item(2, 5).Value = "  +0.43%  "

$ws.Cells.Item(3, 4).Value = "3.482.74"
$ws.Cells.Item(3, 5).Value = "  -0.54%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "593.72"
$ws.Cells.Item(5, 5).Value = "  -0.18%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "178.21"
$ws.Cells.Item(6, 5).Value = "  +3.25%  "

$ws.Cells.Item(8, 5).Value = "  +1.41%  "

$ws.Cells.Item(9, 5).Value = "  -0.40%  "

$ws.Cells.Item(10, 5).Value = "  +4.36%  "

$ws.Cells.Item(11, 5).Value = "  -2.27%  "

$ws.Cells.Item(12, 5).Value = "  +0.49%  "

$ws.Cells.Item(13, 5).Value = "  -0.45%  "

$ws.Cells.Item(14, 5).Value = "  +9.10%  "

$ws.Cells.Item(15, 5).Value = "  +0.87%  "

$ws.Cells.Item(16, 4).Value = "67.307.82"
$ws.Cells.Item(16, 5).Value = "  +0.43%  "

$ws.Cells.Item(17, 5).Value = "  -1.24%  "

$ws.Cells.Item(18, 5).Value = "  -2.33%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.24"
$ws.Cells.Item(19, 5).Value = "  -0.72%  "

$ws.Cells.Item(20, 5).Value = "  +0.04%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "388.42"
$ws.Cells.Item(21, 5).Value = "  -1.94%  "

$ws.Cells.Item(22, 5).Value = "  +0.04%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "73.84"
$ws.Cells.Item(23, 5).Value = "  +0.47%  "

$ws.Cells.Item(24, 5).Value = "  +0.15%  "

$ws.Cells.Item(25, 5).Value = "  +0.08%  "

$ws.Cells.Item(27, 5).Value = "  -0.70%  "

$ws.Cells.Item(28, 5).Value = "  +0.44%  "

$ws.Cells.Item(29, 5).Value = "  -4.09%  "

$ws.Cells.Item(30, 5).Value = "  +0.38%  "

$ws.Cells.Item(31, 5).Value = "  -1.39%  "

$ws.Cells.Item(32, 5).Value = "  -1.05%  "

$ws.Cells.Item(33, 5).Value = "  -0.48%  "

$ws.Cells.Item(34, 5).Value = "  -1.11%  "

$ws.Cells.Item(36, 5).Value = "  +0.04%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.59"
$ws.Cells.Item(37, 5).Value = "  -2.06%  "

$ws.Cells.Item(38, 5).Value = "  +1.10%  "

$ws.Cells.Item(39, 5).Value = "  -1.28%  "

$ws.Cells.Item(40, 2).Value = "dogwifhat"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.75"
$ws.Cells.Item(40, 5).Value = "  +7.79%  "

$ws.Cells.Item(41, 2).Value = "Stacks"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.86"
$ws.Cells.Item(41, 5).Value = "  -2.42%  "

$ws.Cells.Item(42, 5).Value = "  -1.76%  "

$ws.Cells.Item(43, 5).Value = "  -0.19%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "27.03"
$ws.Cells.Item(45, 5).Value = "  -0.55%  "

$ws.Cells.Item(46, 5).Value = "  -0.67%  "

$ws.Cells.Item(47, 5).Value = "  -2.63%  "

$ws.Cells.Item(48, 5).Value = "  -2.83%  "

$ws.Cells.Item(49, 5).Value = "  -0.67%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "334.68"
$ws.Cells.Item(50, 5).Value = "  -1.45%  "

$ws.Cells.Item(51, 5).Value = "  -2.52%  "

